# Fruta / hortaliza, semanal
# Adds a new weekly price report (3 rows) for Granada / "Wonderfull" variety,
# "$/caja 18 kilos granel" unit, "Región de O'Higgins" origin, dated 2022-05-04
# (serial 44685), at the top of the existing data block, pushing the rest of
# the table down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows right above the first data row (row 2), shifting
# all existing data rows down by 3 (old row 2 -> new row 5, etc.)
$ws.Rows("2:4").Insert()

# The inserted rows pick up formatting from the row above (the bold header).
# Reset them to the plain/default formatting used by the rest of the data rows.
$ws.Rows("2:4").ClearFormats()

# Re-apply the date number format used throughout column D.
$ws.Range("D2:D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 2: Wonderfull / Especial
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C2").Value = "Metropolitana"
$ws.Range("D2").Value = 44685
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104001
$ws.Range("J2").Value = "Granada"
$ws.Range("K2").Value = "Wonderfull"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 350
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 21000
$ws.Range("Q2").Value = "`$/caja 18 kilos granel"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1167
$ws.Range("T2").Value = 18

# Row 3: Wonderfull / Primera
$ws.Range("A3").Value = 9
$ws.Range("B3").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C3").Value = "Metropolitana"
$ws.Range("D3").Value = 44685
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100104
$ws.Range("H3").Value = "Frutos de pepita"
$ws.Range("I3").Value = 100104001
$ws.Range("J3").Value = "Granada"
$ws.Range("K3").Value = "Wonderfull"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 330
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "`$/caja 18 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 833
$ws.Range("T3").Value = 18

# Row 4: Wonderfull / Segunda
$ws.Range("A4").Value = 9
$ws.Range("B4").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44685
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104001
$ws.Range("J4").Value = "Granada"
$ws.Range("K4").Value = "Wonderfull"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 280
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("Q4").Value = "`$/caja 18 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 556
$ws.Range("T4").Value = 18
